# Update account to test account: ngqautotester@hpe.com
#
# The sheet had a row (row 2) pointing at a real tester's mailbox
# (yu.li9@hpe.com, with a mailto hyperlink) and a row (row 3) holding the
# generic ngqautotester@hpe.com account (no hyperlink). This swaps things so
# row 2 now holds the plain ngqautotester@hpe.com account/password-hash pair
# (no hyperlink) and row 3 is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlink that lived on F2.
$ws.Hyperlinks.Delete()

# Row 2 becomes the generic test account + its password hash.
$ws.Range("F2").Value = "ngqautotester@hpe.com"
$ws.Range("G2").Value = "58d22393e6aed6fe0e28a9b20e63a85e883b946af9f1ca765172"

# Row 3 (which used to hold the very same test account/hash) is now empty.
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()

# Reflect the last selected cell as saved in the workbook.
$ws.Range("G13").Select()
